$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the import-sample header: "EmpCd" was cryptic -> rename to "Employee Code"
$ws.Range("A1").Value = "Employee Code"

# Update the saved cursor/selection position on the sheet
$ws.Range("E7").Select()
